$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in player names (column A) and positions (column B) for rows 2-13
$players = @(
    @{Row=2;  Name="Felipinho"; Pos="GK"},
    @{Row=3;  Name="Zerefly";   Pos="ZAG"},
    @{Row=4;  Name="Alvestruz"; Pos="ZAG"},
    @{Row=5;  Name="Palomares"; Pos="ZAG"},
    @{Row=6;  Name="Sales";     Pos="ZAG"},
    @{Row=7;  Name="Fabion";    Pos="MC"},
    @{Row=8;  Name="Andrey";    Pos="MC"},
    @{Row=9;  Name="Costa";     Pos="MC"},
    @{Row=10; Name="Viana";     Pos="MC"},
    @{Row=11; Name="Josefino";  Pos="ST"},
    @{Row=12; Name="Pipe";      Pos="ST"},
    @{Row=13; Name="Gio";       Pos="ST"}
)

foreach ($p in $players) {
    $ws.Cells.Item($p.Row, 1).Value = $p.Name
    $ws.Cells.Item($p.Row, 2).Value = $p.Pos
}

# Update the active selection to A13
$ws.Range("A13").Select()
